$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header timestamp
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 04:33"

# Panama overtakes Sudafrica (rows 49/50 swap, Panama gets new figures)
$ws.Range("A49").Value = "Panama"
$ws.Range("B49").Value = 7387
$ws.Range("C49").Value = 190
$ws.Range("D49").Value = 726
$ws.Range("E49").Value = 6458
$ws.Range("F49").Value = 93
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 203

$ws.Range("A50").Value = "Sudafrica"
$ws.Range("B50").Value = 7220
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 2746
$ws.Range("E50").Value = 4336
$ws.Range("F50").Value = 36
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 138

# Australia data refresh (no reorder)
$ws.Range("B51").Value = 6847
$ws.Range("C51").Value = 22
$ws.Range("D51").Value = 5886
$ws.Range("E51").Value = 866
$ws.Range("F51").Value = 28
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 95

# San Cristobal y Nieves overtakes Burundi (rows 198/199 swap)
$ws.Range("A198").Value = "San Cristobal y Nieves"
$ws.Range("B198").Value = 15
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 8
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Burundi"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 7
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

# Mauritania overtakes Papua Nueva Guinea (rows 209/210 swap, Papua gets new figures)
$ws.Range("A209").Value = "Mauritania"
$ws.Range("B209").Value = 8
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 6
$ws.Range("E209").Value = 1
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Papua Nueva Guinea"
$ws.Range("B210").Value = 8
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 8
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
